# Update ABA_Questionnaire Recommendation Engine Improvement.xlsx
# - rewrite a few free-text "activity" answers (col G) on the
#   "Questionnaire Results" sheet
# - add two new participants (Tino, Lou) as rows 12 & 13
# - move the active selection to A14

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Questionnaire Results")

$xlCenter = -4108

# --- Rewrite a few existing free-text activity answers (column G) ---
# (order matters: new shared-string entries are appended in first-use
#  order, and this matches the order the author entered them in)
$ws.Range("G8").Value = "Clubbing"
$ws.Range("G11").Value = "Outdoor activity like hiking"
$ws.Range("G5").Value = "Shopping"

# --- Add participant 9: Tino ---
$ws.Range("A12").Value = 9
$ws.Range("A12").HorizontalAlignment = $xlCenter
$ws.Range("B12").Value = "Tino"
$ws.Range("C12").Value = "B"
$ws.Range("C12").HorizontalAlignment = $xlCenter
$ws.Range("D12").Value = "B"
$ws.Range("D12").HorizontalAlignment = $xlCenter
$ws.Range("E12").Value = "C"
$ws.Range("E12").HorizontalAlignment = $xlCenter
$ws.Range("F12").Value = "Tipp A"
$ws.Range("F12").HorizontalAlignment = $xlCenter
$ws.Range("G12").Value = "Shopping"
$ws.Range("G12").HorizontalAlignment = $xlCenter
$ws.Range("I12").Value = "Tipp A"
$ws.Range("I12").HorizontalAlignment = $xlCenter
# nudge the shared "Match" formula so it re-evaluates against the new row
$ws.Range("J12").Formula = $ws.Range("J12").Formula

# --- Add participant 10: Lou ---
$ws.Range("A13").Value = 10
$ws.Range("A13").HorizontalAlignment = $xlCenter
$ws.Range("B13").Value = "Lou"
$ws.Range("C13").Value = "B"
$ws.Range("C13").HorizontalAlignment = $xlCenter
$ws.Range("D13").Value = "A"
$ws.Range("D13").HorizontalAlignment = $xlCenter
$ws.Range("E13").Value = "A"
$ws.Range("E13").HorizontalAlignment = $xlCenter
$ws.Range("F13").Value = "Tipp A"
$ws.Range("F13").HorizontalAlignment = $xlCenter
$ws.Range("G13").Value = "Shopping"
$ws.Range("G13").HorizontalAlignment = $xlCenter
$ws.Range("I13").Value = "Tipp B"
$ws.Range("I13").HorizontalAlignment = $xlCenter
# nudge the shared "Match" formula so it re-evaluates against the new row
$ws.Range("J13").Formula = $ws.Range("J13").Formula

# --- Move the selection like the author left it ---
$ws.Range("A14").Select()
